# Ahora funciona, hay que cambiar la clase Ui_Caja
# Adds one new transaction row to each sheet involved in a "Caja" (checkout) operation.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, [string]$text)
    # Force a numeric-looking string to be stored as text instead of being
    # auto-coerced to a number, without leaving a residual custom style on
    # the cell once we are done.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---- Clientes (sheet1) : new row 11 ----
$ws = $wb.Worksheets.Item("Clientes")
$ws.Cells.Item(11, 1).Value = 1234
$ws.Cells.Item(11, 2).Value = "migue"
Set-TextValue $ws.Cells.Item(11, 3) "123456789"

# ---- Productos (sheet2) : new row 10 ----
$ws = $wb.Worksheets.Item("Productos")
$ws.Cells.Item(10, 1).Value = "REF123"
Set-TextValue $ws.Cells.Item(10, 2) "1234567890123"
$ws.Cells.Item(10, 3).Value = "Marca A"
$ws.Cells.Item(10, 4).Value = 10
$ws.Cells.Item(10, 5).Value = 20
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = $false
$ws.Cells.Item(10, 8).Value = "31/05/2024 23:09"

# ---- VentaProductos (sheet4) : new row 10 ----
$ws = $wb.Worksheets.Item("VentaProductos")
$ws.Cells.Item(10, 1).Value = 123
$ws.Cells.Item(10, 2).Value = 1234
$ws.Cells.Item(10, 3).Value = "migue"
$ws.Cells.Item(10, 4).Value = "Shampoo"
$ws.Cells.Item(10, 5).Value = "31/05/2024 23:09"
$ws.Cells.Item(10, 6).Value = 12
$ws.Cells.Item(10, 7).Value = 3221
$ws.Cells.Item(10, 8).Value = "efectivo"

# ---- ReservasServicios (sheet7) : new row 10 ----
$ws = $wb.Worksheets.Item("ReservasServicios")
$ws.Cells.Item(10, 1).Value = 12
$ws.Cells.Item(10, 2).Value = "mgiue"
$ws.Cells.Item(10, 3).Value = 345
$ws.Cells.Item(10, 4).Value = "31/05/2024 23:09"
$ws.Cells.Item(10, 5).Value = "28/05/2024 14:30"

# ---- Facturas (sheet8) : new row 10 ----
$ws = $wb.Worksheets.Item("Facturas")
$ws.Cells.Item(10, 1).Value = 1234
$ws.Cells.Item(10, 2).Value = "migue"
$ws.Cells.Item(10, 3).Value = "31/05/2024 23:09"
$ws.Cells.Item(10, 4).Value = 3221
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 3221
$ws.Cells.Item(10, 7).Value = 1

# ---- Usuarios (sheet9) : new row 11 ----
$ws = $wb.Worksheets.Item("Usuarios")
$ws.Cells.Item(11, 1).Value = 2
$ws.Cells.Item(11, 2).Value = "admin"
$ws.Cells.Item(11, 3).Value = 12345
$ws.Cells.Item(11, 4).Value = 2
